# Auto-generated edit script: updates commodity price/profit figures
# across the Pandaemonium_Profits sheets (scheduled price-refresh run).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Cells.Item(19, 8).Value = 22223262
$ws.Cells.Item(19, 9).Value = 166667020
$ws.Cells.Item(19, 10).Value = 1144
$ws.Cells.Item(19, 11).Value = 166667020
$ws.Cells.Item(19, 12).Value = 1144
$ws.Cells.Item(19, 13).Value = -166666845
$ws.Cells.Item(19, 14).Value = -1494

# Row 53
$ws.Cells.Item(53, 8).Value = 383.1
$ws.Cells.Item(53, 9).Value = 358.125
$ws.Cells.Item(53, 10).Value = 483
$ws.Cells.Item(53, 11).Value = 358.125
$ws.Cells.Item(53, 12).Value = 483
$ws.Cells.Item(53, 13).Value = 278.875
$ws.Cells.Item(53, 14).Value = -1757

# Row 55
$ws.Cells.Item(55, 8).Value = 79.09999999999999
$ws.Cells.Item(55, 9).Value = 80.111115
$ws.Cells.Item(55, 10).Value = 70
$ws.Cells.Item(55, 11).Value = 80.111115
$ws.Cells.Item(55, 12).Value = 70
$ws.Cells.Item(55, 13).Value = 133.888885
$ws.Cells.Item(55, 14).Value = -498

# Row 64
$ws.Cells.Item(64, 8).Value = 4015.6
$ws.Cells.Item(64, 9).Value = 3593.125
$ws.Cells.Item(64, 10).Value = 4766.6665
$ws.Cells.Item(64, 11).Value = 3593.125
$ws.Cells.Item(64, 12).Value = 4766.6665
$ws.Cells.Item(64, 13).Value = -3345.125
$ws.Cells.Item(64, 14).Value = -5262.6665

# Row 67
$ws.Cells.Item(67, 8).Value = 4015.6
$ws.Cells.Item(67, 9).Value = 3593.125
$ws.Cells.Item(67, 10).Value = 4766.6665
$ws.Cells.Item(67, 11).Value = 3593.125
$ws.Cells.Item(67, 12).Value = 4766.6665
$ws.Cells.Item(67, 13).Value = -2735.125
$ws.Cells.Item(67, 14).Value = -6482.6665

# Row 98
$ws.Cells.Item(98, 8).Value = 4632.778
$ws.Cells.Item(98, 9).Value = 3099.2856
$ws.Cells.Item(98, 10).Value = 10000
$ws.Cells.Item(98, 11).Value = 3099.2856
$ws.Cells.Item(98, 12).Value = 10000
$ws.Cells.Item(98, 13).Value = -1601.2856
$ws.Cells.Item(98, 14).Value = -12996

# Row 122
$ws.Cells.Item(122, 8).Value = 4632.778
$ws.Cells.Item(122, 9).Value = 3099.2856
$ws.Cells.Item(122, 10).Value = 10000
$ws.Cells.Item(122, 11).Value = 9297.856800000001
$ws.Cells.Item(122, 12).Value = 30000
$ws.Cells.Item(122, 13).Value = -6847.856800000001
$ws.Cells.Item(122, 14).Value = -34900

# Row 132
$ws.Cells.Item(132, 8).Value = 1546.3036
$ws.Cells.Item(132, 9).Value = 1304.0209
$ws.Cells.Item(132, 10).Value = 3000
$ws.Cells.Item(132, 11).Value = 3912.0627
$ws.Cells.Item(132, 12).Value = 9000
$ws.Cells.Item(132, 13).Value = -1382.0627

# Row 137
$ws.Cells.Item(137, 8).Value = 3190.3555
$ws.Cells.Item(137, 9).Value = 1560.2069
$ws.Cells.Item(137, 10).Value = 6145
$ws.Cells.Item(137, 11).Value = 4680.620699999999
$ws.Cells.Item(137, 12).Value = 18435
$ws.Cells.Item(137, 13).Value = -2130.620699999999
$ws.Cells.Item(137, 14).Value = -23535

# Row 140
$ws.Cells.Item(140, 8).Value = 82509.336
$ws.Cells.Item(140, 9).Value = 0
$ws.Cells.Item(140, 10).Value = 82509.336
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(140, 12).Value = 82509.336
$ws.Cells.Item(140, 14).Value = -92869.336

# Row 141
$ws.Cells.Item(141, 8).Value = 3645
$ws.Cells.Item(141, 9).Value = 3699.1667
$ws.Cells.Item(141, 10).Value = 3598.5715
$ws.Cells.Item(141, 11).Value = 11097.5001
$ws.Cells.Item(141, 12).Value = 10795.7145
$ws.Cells.Item(141, 13).Value = -5917.500100000001
$ws.Cells.Item(141, 14).Value = -21155.7145

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 1539.081
$ws.Cells.Item(2, 9).Value = 1551.5333
$ws.Cells.Item(2, 10).Value = 1485.7142
$ws.Cells.Item(2, 11).Value = 1551.5333
$ws.Cells.Item(2, 12).Value = 1485.7142
$ws.Cells.Item(2, 13).Value = -1438.5333

# Row 32
$ws.Cells.Item(32, 8).Value = 7407.625
$ws.Cells.Item(32, 9).Value = 7087.7637
$ws.Cells.Item(32, 10).Value = 25000
$ws.Cells.Item(32, 11).Value = 7087.7637
$ws.Cells.Item(32, 12).Value = 25000
$ws.Cells.Item(32, 13).Value = -6800.7637

# Row 74
$ws.Cells.Item(74, 8).Value = 103136.164
$ws.Cells.Item(74, 9).Value = 116173.805
$ws.Cells.Item(74, 10).Value = 22302.8
$ws.Cells.Item(74, 11).Value = 116173.805
$ws.Cells.Item(74, 12).Value = 22302.8
$ws.Cells.Item(74, 13).Value = -115299.805
$ws.Cells.Item(74, 14).Value = -24050.8

# Row 77
$ws.Cells.Item(77, 8).Value = 103136.164
$ws.Cells.Item(77, 9).Value = 116173.805
$ws.Cells.Item(77, 10).Value = 22302.8
$ws.Cells.Item(77, 11).Value = 580869.0249999999
$ws.Cells.Item(77, 12).Value = 111514
$ws.Cells.Item(77, 13).Value = -576501.0249999999
$ws.Cells.Item(77, 14).Value = -120250

# Row 102
$ws.Cells.Item(102, 8).Value = 2540
$ws.Cells.Item(102, 9).Value = 2031.1111
$ws.Cells.Item(102, 10).Value = 4066.6667
$ws.Cells.Item(102, 11).Value = 2031.1111
$ws.Cells.Item(102, 12).Value = 4066.6667
$ws.Cells.Item(102, 13).Value = -409.1111000000001

# Row 116
$ws.Cells.Item(116, 8).Value = 1539.081
$ws.Cells.Item(116, 9).Value = 1551.5333
$ws.Cells.Item(116, 10).Value = 1485.7142
$ws.Cells.Item(116, 11).Value = 1551.5333
$ws.Cells.Item(116, 12).Value = 1485.7142
$ws.Cells.Item(116, 13).Value = 742.4666999999999

# Row 132
$ws.Cells.Item(132, 8).Value = 4651.449
$ws.Cells.Item(132, 9).Value = 1613
$ws.Cells.Item(132, 10).Value = 7816.5
$ws.Cells.Item(132, 11).Value = 4839
$ws.Cells.Item(132, 12).Value = 23449.5
$ws.Cells.Item(132, 13).Value = -2309

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 1539.081
$ws.Cells.Item(3, 9).Value = 1551.5333
$ws.Cells.Item(3, 10).Value = 1485.7142
$ws.Cells.Item(3, 11).Value = 1551.5333
$ws.Cells.Item(3, 12).Value = 1485.7142
$ws.Cells.Item(3, 13).Value = -1437.5333

# Row 94
$ws.Cells.Item(94, 8).Value = 1660.2894
$ws.Cells.Item(94, 9).Value = 1572.4231
$ws.Cells.Item(94, 10).Value = 1850.6666
$ws.Cells.Item(94, 11).Value = 1572.4231
$ws.Cells.Item(94, 12).Value = 1850.6666
$ws.Cells.Item(94, 13).Value = -1121.4231
$ws.Cells.Item(94, 14).Value = -2752.6666

# Row 107
$ws.Cells.Item(107, 8).Value = 2621
$ws.Cells.Item(107, 9).Value = 2415.8823
$ws.Cells.Item(107, 10).Value = 3202.1667
$ws.Cells.Item(107, 11).Value = 2415.8823
$ws.Cells.Item(107, 12).Value = 3202.1667
$ws.Cells.Item(107, 13).Value = -495.8823000000002

# Row 134
$ws.Cells.Item(134, 8).Value = 32587.666
$ws.Cells.Item(134, 9).Value = 2398.6155
$ws.Cells.Item(134, 10).Value = 144718.42
$ws.Cells.Item(134, 11).Value = 7195.8465
$ws.Cells.Item(134, 12).Value = 434155.26
$ws.Cells.Item(134, 13).Value = -4660.8465
$ws.Cells.Item(134, 14).Value = -439225.26

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Cells.Item(7, 8).Value = 214
$ws.Cells.Item(7, 9).Value = 199.75
$ws.Cells.Item(7, 10).Value = 233
$ws.Cells.Item(7, 11).Value = 199.75
$ws.Cells.Item(7, 12).Value = 233
$ws.Cells.Item(7, 13).Value = -86.75
$ws.Cells.Item(7, 14).Value = -459

# Row 31
$ws.Cells.Item(31, 8).Value = 2652.5
$ws.Cells.Item(31, 9).Value = 2079.5588
$ws.Cells.Item(31, 10).Value = 3261.25
$ws.Cells.Item(31, 11).Value = 2079.5588
$ws.Cells.Item(31, 12).Value = 3261.25
$ws.Cells.Item(31, 13).Value = -1784.5588
$ws.Cells.Item(31, 14).Value = -3851.25

# Row 34
$ws.Cells.Item(34, 8).Value = 2652.5
$ws.Cells.Item(34, 9).Value = 2079.5588
$ws.Cells.Item(34, 10).Value = 3261.25
$ws.Cells.Item(34, 11).Value = 2079.5588
$ws.Cells.Item(34, 12).Value = 3261.25
$ws.Cells.Item(34, 13).Value = -1877.5588
$ws.Cells.Item(34, 14).Value = -3665.25

# Row 58
$ws.Cells.Item(58, 8).Value = 2528359
$ws.Cells.Item(58, 9).Value = 3954979.2
$ws.Cells.Item(58, 10).Value = 4338.4614
$ws.Cells.Item(58, 11).Value = 3954979.2
$ws.Cells.Item(58, 12).Value = 4338.4614
$ws.Cells.Item(58, 13).Value = -3954776.2
$ws.Cells.Item(58, 14).Value = -4744.4614

# Row 62
$ws.Cells.Item(62, 8).Value = 3167.889
$ws.Cells.Item(62, 9).Value = 3000.8333
$ws.Cells.Item(62, 10).Value = 3502
$ws.Cells.Item(62, 11).Value = 3000.8333
$ws.Cells.Item(62, 12).Value = 3502
$ws.Cells.Item(62, 13).Value = -2376.8333
$ws.Cells.Item(62, 14).Value = -4750

# Row 65
$ws.Cells.Item(65, 8).Value = 3167.889
$ws.Cells.Item(65, 9).Value = 3000.8333
$ws.Cells.Item(65, 10).Value = 3502
$ws.Cells.Item(65, 11).Value = 15004.1665
$ws.Cells.Item(65, 12).Value = 17510
$ws.Cells.Item(65, 13).Value = -11884.1665
$ws.Cells.Item(65, 14).Value = -23750

# Row 94
$ws.Cells.Item(94, 8).Value = 1472.9166
$ws.Cells.Item(94, 9).Value = 1506.7142
$ws.Cells.Item(94, 10).Value = 1425.6
$ws.Cells.Item(94, 11).Value = 1506.7142
$ws.Cells.Item(94, 12).Value = 1425.6
$ws.Cells.Item(94, 13).Value = -1055.7142
$ws.Cells.Item(94, 14).Value = -2327.6

# Row 122
$ws.Cells.Item(122, 8).Value = 8097.143
$ws.Cells.Item(122, 9).Value = 4425.846
$ws.Cells.Item(122, 10).Value = 11278.934
$ws.Cells.Item(122, 11).Value = 13277.538
$ws.Cells.Item(122, 12).Value = 33836.802
$ws.Cells.Item(122, 13).Value = -10827.538
$ws.Cells.Item(122, 14).Value = -38736.802

# Row 132
$ws.Cells.Item(132, 8).Value = 2714.2727
$ws.Cells.Item(132, 9).Value = 2296.4119
$ws.Cells.Item(132, 10).Value = 4135
$ws.Cells.Item(132, 11).Value = 6889.2357
$ws.Cells.Item(132, 12).Value = 12405
$ws.Cells.Item(132, 13).Value = -4359.2357
$ws.Cells.Item(132, 14).Value = -17465

# Row 134
$ws.Cells.Item(134, 8).Value = 23395.6
$ws.Cells.Item(134, 9).Value = 54746.367
$ws.Cells.Item(134, 10).Value = 4180.613
$ws.Cells.Item(134, 11).Value = 164239.101
$ws.Cells.Item(134, 12).Value = 12541.839
$ws.Cells.Item(134, 13).Value = -161704.101
$ws.Cells.Item(134, 14).Value = -17611.839

# Row 136
$ws.Cells.Item(136, 8).Value = 2528359
$ws.Cells.Item(136, 9).Value = 3954979.2
$ws.Cells.Item(136, 10).Value = 4338.4614
$ws.Cells.Item(136, 11).Value = 11864937.6
$ws.Cells.Item(136, 12).Value = 13015.3842
$ws.Cells.Item(136, 13).Value = -11862387.6
$ws.Cells.Item(136, 14).Value = -18115.3842

# Row 140
$ws.Cells.Item(140, 8).Value = 38853.332
$ws.Cells.Item(140, 9).Value = 0
$ws.Cells.Item(140, 10).Value = 38853.332
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(140, 12).Value = 38853.332
$ws.Cells.Item(140, 13).ClearContents()
$ws.Cells.Item(140, 14).Value = -49213.332

$ws = $wb.Worksheets.Item("CUL")
# Row 13
$ws.Cells.Item(13, 8).Value = 300
$ws.Cells.Item(13, 9).Value = 300
$ws.Cells.Item(13, 10).Value = 300
$ws.Cells.Item(13, 11).Value = 900
$ws.Cells.Item(13, 12).Value = 900
$ws.Cells.Item(13, 13).Value = -732
$ws.Cells.Item(13, 14).Value = -1236

# Row 92
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 13).ClearContents()
$ws.Cells.Item(92, 14).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Cells.Item(2, 8).Value = 294.66666
$ws.Cells.Item(2, 9).Value = 122.666664
$ws.Cells.Item(2, 10).Value = 466.66666
$ws.Cells.Item(2, 11).Value = 122.666664
$ws.Cells.Item(2, 12).Value = 466.66666
$ws.Cells.Item(2, 13).Value = -9.666663999999997
$ws.Cells.Item(2, 14).Value = -692.66666

# Row 97
$ws.Cells.Item(97, 8).Value = 1499.6451
$ws.Cells.Item(97, 9).Value = 1194.5416
$ws.Cells.Item(97, 10).Value = 2545.7144
$ws.Cells.Item(97, 11).Value = 1194.5416
$ws.Cells.Item(97, 12).Value = 2545.7144
$ws.Cells.Item(97, 13).Value = -698.5416
$ws.Cells.Item(97, 14).Value = -3537.7144

# Row 126
$ws.Cells.Item(126, 8).Value = 3040.45
$ws.Cells.Item(126, 9).Value = 1989
$ws.Cells.Item(126, 10).Value = 3741.4167
$ws.Cells.Item(126, 11).Value = 5967
$ws.Cells.Item(126, 12).Value = 11224.2501
$ws.Cells.Item(126, 13).Value = -3497
$ws.Cells.Item(126, 14).Value = -16164.2501

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Cells.Item(40, 8).Value = 4523.357
$ws.Cells.Item(40, 9).Value = 4049.111
$ws.Cells.Item(40, 10).Value = 5377
$ws.Cells.Item(40, 11).Value = 4049.111
$ws.Cells.Item(40, 12).Value = 5377
$ws.Cells.Item(40, 13).Value = -3913.111
$ws.Cells.Item(40, 14).Value = -5649

# Row 93
$ws.Cells.Item(93, 8).Value = 2379.75
$ws.Cells.Item(93, 9).Value = 2314.3845
$ws.Cells.Item(93, 10).Value = 2501.1428
$ws.Cells.Item(93, 11).Value = 2314.3845
$ws.Cells.Item(93, 12).Value = 2501.1428
$ws.Cells.Item(93, 13).Value = -1066.3845
$ws.Cells.Item(93, 14).Value = -4997.1428

# Row 100
$ws.Cells.Item(100, 8).Value = 4509.5264
$ws.Cells.Item(100, 9).Value = 2598.4167
$ws.Cells.Item(100, 10).Value = 7785.7144
$ws.Cells.Item(100, 11).Value = 2598.4167
$ws.Cells.Item(100, 12).Value = 7785.7144
$ws.Cells.Item(100, 13).Value = -2057.4167
$ws.Cells.Item(100, 14).Value = -8867.714400000001

# Row 136
$ws.Cells.Item(136, 8).Value = 4341.5576
$ws.Cells.Item(136, 9).Value = 2689.3235
$ws.Cells.Item(136, 10).Value = 7462.4443
$ws.Cells.Item(136, 11).Value = 8067.970499999999
$ws.Cells.Item(136, 12).Value = 22387.3329
$ws.Cells.Item(136, 13).Value = -5517.970499999999

# Row 140
$ws.Cells.Item(140, 8).Value = 73150.78
$ws.Cells.Item(140, 9).Value = 0
$ws.Cells.Item(140, 10).Value = 73150.78
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(140, 12).Value = 73150.78
$ws.Cells.Item(140, 14).Value = -83510.78

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Cells.Item(132, 8).Value = 3493.3142
$ws.Cells.Item(132, 9).Value = 3149.65
$ws.Cells.Item(132, 10).Value = 3951.5334
$ws.Cells.Item(132, 11).Value = 9448.950000000001
$ws.Cells.Item(132, 12).Value = 11854.6002
$ws.Cells.Item(132, 13).Value = -6918.950000000001
$ws.Cells.Item(132, 14).Value = -16914.6002

# Row 136
$ws.Cells.Item(136, 8).Value = 4480.113
$ws.Cells.Item(136, 9).Value = 1923.5883
$ws.Cells.Item(136, 10).Value = 7584.4644
$ws.Cells.Item(136, 11).Value = 5770.7649
$ws.Cells.Item(136, 12).Value = 22753.3932
$ws.Cells.Item(136, 13).Value = -3220.7649
$ws.Cells.Item(136, 14).Value = -27853.3932

